$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''27.023.57'
$ws.Range('E2').Value = '''  -1.04%  '
$ws.Range('D3').Value = '''1.821.26'
$ws.Range('E3').Value = '''  -0.26%  '
$ws.Range('D4').Value = '''1.002'
$ws.Range('E4').Value = '''  -0.41%  '
$ws.Range('D5').Value = '''310.75'
$ws.Range('E5').Value = '''  -1.20%  '
$ws.Range('E6').Value = '''  -0.37%  '
$ws.Range('D7').Value = '''0.4477'
$ws.Range('E7').Value = '''  +4.77%  '
$ws.Range('D8').Value = '''0.3696'
$ws.Range('E8').Value = '''  +0.12%  '
$ws.Range('D9').Value = '''0.07304'
$ws.Range('E9').Value = '''  +0.70%  '
$ws.Range('D10').Value = '''0.8563'
$ws.Range('E10').Value = '''  -0.96%  '
$ws.Range('D11').Value = '''20.74'
$ws.Range('E11').Value = '''  -1.59%  '
$ws.Range('D12').Value = '''1.825.44'
$ws.Range('E12').Value = '''  -0.11%  '
$ws.Range('D13').Value = '''6.639'
$ws.Range('E13').Value = '''  -1.07%  '
$ws.Range('B14').Value = 'TRON'
$ws.Range('C14').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D14').Value = '''0.07104'
$ws.Range('E14').Value = '''  +0.02%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = '''5.329'
$ws.Range('E15').Value = '''  +0.23%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').Value = '''92.14'
$ws.Range('E16').Value = '''  +3.97%  '
$ws.Range('E17').Value = '''  -0.42%  '
$ws.Range('D18').Value = '''0.000008782'
$ws.Range('E18').Value = '''  -0.98%  '
$ws.Range('E19').Value = '''  -0.47%  '
$ws.Range('D20').Value = '''14.96'
$ws.Range('E20').Value = '''  -0.92%  '
$ws.Range('D21').Value = '''27.051.81'
$ws.Range('E21').Value = '''  -1.05%  '
$ws.Range('D22').Value = '''5.166'
$ws.Range('E22').Value = '''  +0.41%  '
$ws.Range('D23').Value = '''10.91'
$ws.Range('E23').Value = '''  +0.20%  '
$ws.Range('D24').Value = '''1.988'
$ws.Range('E24').Value = '''  -0.83%  '
$ws.Range('D25').Value = '''151.74'
$ws.Range('E25').Value = '''  -0.91%  '
$ws.Range('D26').Value = '''2.228'
$ws.Range('E26').Value = '''  +3.71%  '
$ws.Range('D27').Value = '''18.48'
$ws.Range('E27').Value = '''  +0.69%  '
$ws.Range('D28').Value = '''5.248'
$ws.Range('E28').Value = '''  -0.11%  '
$ws.Range('D29').Value = '''116.62'
$ws.Range('E29').Value = '''  +0.02%  '
$ws.Range('D30').Value = '''0.08846'
$ws.Range('E30').Value = '''  -0.65%  '
$ws.Range('D31').Value = '''0.7535'
$ws.Range('E31').Value = '''  -0.49%  '
$ws.Range('D32').Value = '''1.181'
$ws.Range('E32').Value = '''  -1.73%  '
$ws.Range('D33').Value = '''2.962'
$ws.Range('E33').Value = '''  +4.34%  '
$ws.Range('D34').Value = '''4.452'
$ws.Range('E34').Value = '''  -0.16%  '
$ws.Range('D35').Value = '''1.001'
$ws.Range('E35').Value = '''  -0.44%  '
$ws.Range('E36').Value = '''  -1.76%  '
$ws.Range('D37').Value = '''0.01966'
$ws.Range('E37').Value = '''  -0.71%  '
$ws.Range('D38').Value = '''0.05235'
$ws.Range('E38').Value = '''  -0.80%  '
$ws.Range('D39').Value = '''0.5311'
$ws.Range('E39').Value = '''  +5.12%  '
$ws.Range('D40').Value = '''2.886'
$ws.Range('E40').Value = '''  +0.37%  '
$ws.Range('D41').Value = '''7.126'
$ws.Range('E41').Value = '''  -0.57%  '
$ws.Range('E42').Value = '''  +0.52%  '
$ws.Range('D43').Value = '''0.5232'
$ws.Range('E43').Value = '''  +10.00%  '
$ws.Range('D44').Value = '''8.499'
$ws.Range('E44').Value = '''  -2.05%  '
$ws.Range('D45').Value = '''10.61'
$ws.Range('E45').Value = '''  +0.14%  '
$ws.Range('D46').Value = '''1.969'
$ws.Range('E46').Value = '''  +7.46%  '
$ws.Range('D47').Value = '''105.48'
$ws.Range('E47').Value = '''  -2.09%  '
$ws.Range('E48').Value = '''  -0.47%  '
$ws.Range('D49').Value = '''1.668'
$ws.Range('E49').Value = '''  +0.07%  '
$ws.Range('D50').Value = '''0.06388'
$ws.Range('E50').Value = '''  +0.21%  '
$ws.Range('D51').Value = '''0.9193'
$ws.Range('E51').Value = '''  +0.10%  '
